$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings in column D remain text (matching source format)
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '30.623.64'
$ws.Range("E2").Value = '  +0.95%  '
$ws.Range("D3").Value = '1.873.27'
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = '247.79'
$ws.Range("E5").Value = '  +1.07%  '
$ws.Range("D6").Value = '1.000'
$ws.Range("E6").Value = '  +0.01%  '
$ws.Range("E7").Value = '  +0.09%  '
$ws.Range("D8").Value = '0.2911'
$ws.Range("E8").Value = '  +1.19%  '
$ws.Range("D9").Value = '0.06483'
$ws.Range("E9").Value = '  +0.18%  '
$ws.Range("D10").Value = '22.08'
$ws.Range("E10").Value = '  +4.65%  '
$ws.Range("D11").Value = '0.07700'
$ws.Range("E11").Value = '  -0.96%  '
$ws.Range("B12").Value = 'Polygon'
$ws.Range("C12").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D12").Value = '0.7381'
$ws.Range("E12").Value = '  +1.75%  '
$ws.Range("B13").Value = 'Litecoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D13").Value = '96.53'
$ws.Range("E13").Value = '  +1.44%  '
$ws.Range("D14").Value = '1.870.02'
$ws.Range("E14").Value = '  -0.03%  '
$ws.Range("D15").Value = '5.167'
$ws.Range("E15").Value = '  +0.46%  '
$ws.Range("D16").Value = '273.32'
$ws.Range("E16").Value = '  -0.80%  '
$ws.Range("D17").Value = '30.639.31'
$ws.Range("E17").Value = '  +1.04%  '
$ws.Range("D18").Value = '13.33'
$ws.Range("E18").Value = '  -0.30%  '
$ws.Range("E19").Value = '  +0.02%  '
$ws.Range("D20").Value = '0.000007523'
$ws.Range("E20").Value = '  -0.45%  '
$ws.Range("D21").Value = '2.117.74'
$ws.Range("E21").Value = '  +0.28%  '
$ws.Range("D22").Value = '0.9999'
$ws.Range("E22").Value = '  -0.01%  '
$ws.Range("D23").Value = '5.272'
$ws.Range("E23").Value = '  +0.45%  '
$ws.Range("D24").Value = '6.199'
$ws.Range("E24").Value = '  +0.34%  '
$ws.Range("D25").Value = '9.233'
$ws.Range("E25").Value = '  -0.33%  '
$ws.Range("D26").Value = '164.41'
$ws.Range("E26").Value = '  -0.73%  '
$ws.Range("D27").Value = '18.76'
$ws.Range("E27").Value = '  -0.78%  '
$ws.Range("E28").Value = '  -0.28%  '
$ws.Range("D29").Value = '0.1003'
$ws.Range("E29").Value = '  +1.20%  '
$ws.Range("E30").Value = '  -2.48%  '
$ws.Range("D31").Value = '1.510'
$ws.Range("E31").Value = '  -0.86%  '
$ws.Range("D32").Value = '4.286'
$ws.Range("E32").Value = '  -0.30%  '
$ws.Range("D33").Value = '4.102'
$ws.Range("E33").Value = '  +1.64%  '
$ws.Range("D34").Value = '0.04803'
$ws.Range("E34").Value = '  +0.63%  '
$ws.Range("D35").Value = '1.122'
$ws.Range("E35").Value = '  -0.13%  '
$ws.Range("D36").Value = '0.6963'
$ws.Range("E36").Value = '  -0.26%  '
$ws.Range("D37").Value = '2.719'
$ws.Range("E37").Value = '  -0.07%  '
$ws.Range("D38").Value = '0.01854'
$ws.Range("E38").Value = '  +0.45%  '
$ws.Range("D39").Value = '2.755'
$ws.Range("E39").Value = '  +0.47%  '
$ws.Range("D40").Value = '6.256'
$ws.Range("D41").Value = '73.30'
$ws.Range("E41").Value = '  +4.23%  '
$ws.Range("E42").Value = '  +3.07%  '
$ws.Range("D43").Value = '0.4183'
$ws.Range("E43").Value = '  +1.39%  '
$ws.Range("D44").Value = '1.001'
$ws.Range("E44").Value = '  +0.05%  '
$ws.Range("D45").Value = '0.8351'
$ws.Range("E45").Value = '  -0.64%  '
$ws.Range("D46").Value = '101.97'
$ws.Range("E46").Value = '  -0.46%  '
$ws.Range("D47").Value = '9.399'
$ws.Range("E47").Value = '  +0.23%  '
$ws.Range("D48").Value = '35.49'
$ws.Range("E48").Value = '  +0.66%  '
$ws.Range("D49").Value = '6.990'
$ws.Range("E49").Value = '  -1.63%  '
$ws.Range("D50").Value = '919.92'
$ws.Range("E50").Value = '  +0.01%  '
$ws.Range("E51").Value = '  +1.43%  '
